$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3039.2
$ws.Range("I29").Value = 2099.75
$ws.Range("J29").Value = 3665.5
$ws.Range("K29").Value = 6299.25
$ws.Range("L29").Value = 10996.5
$ws.Range("M29").Value = -6018.25
$ws.Range("N29").Value = -11558.5
$ws.Range("H55").Value = 265.5
$ws.Range("I55").Value = 168.5
$ws.Range("J55").Value = 362.5
$ws.Range("K55").Value = 168.5
$ws.Range("L55").Value = 362.5
$ws.Range("M55").Value = 45.5
$ws.Range("N55").Value = -790.5
$ws.Range("H58").Value = 603.5714
$ws.Range("I58").Value = 603.5714
$ws.Range("K58").Value = 1810.7142
$ws.Range("M58").Value = -1660.7142
$ws.Range("H87").Value = 82611.11
$ws.Range("J87").Value = 82611.11
$ws.Range("L87").Value = 82611.11
$ws.Range("N87").Value = -85107.11
$ws.Range("H90").Value = 82611.11
$ws.Range("J90").Value = 82611.11
$ws.Range("L90").Value = 247833.33
$ws.Range("N90").Value = -260313.33
$ws.Range("H132").Value = 124312.9
$ws.Range("I132").Value = 335717
$ws.Range("J132").Value = 14545.385
$ws.Range("K132").Value = 1007151
$ws.Range("L132").Value = 43636.155
$ws.Range("M132").Value = -1004621
$ws.Range("N132").Value = -48696.155
$ws.Range("H135").Value = 6244.45
$ws.Range("J135").Value = 12232.223
$ws.Range("L135").Value = 110090.007
$ws.Range("N135").Value = -115160.007
$ws.Range("H137").Value = 17548302
$ws.Range("I137").Value = 1767
$ws.Range("K137").Value = 5301
$ws.Range("M137").Value = -2751
$ws.Range("H138").Value = 6627.5415
$ws.Range("I138").Value = 2820.2856
$ws.Range("K138").Value = 8460.856800000001
$ws.Range("M138").Value = -3320.856800000001
$ws.Range("H140").Value = 71656.5
$ws.Range("J140").Value = 69141.25
$ws.Range("L140").Value = 69141.25
$ws.Range("N140").Value = -79501.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 791034.5
$ws.Range("I2").Value = 1086778.9
$ws.Range("J2").Value = 72798.14
$ws.Range("K2").Value = 1086778.9
$ws.Range("L2").Value = 72798.14
$ws.Range("M2").Value = -1086665.9
$ws.Range("N2").Value = -73024.14
$ws.Range("H61").Value = 10647150
$ws.Range("I61").Value = 13164597
$ws.Range("J61").Value = 17930.111
$ws.Range("K61").Value = 13164597
$ws.Range("L61").Value = 17930.111
$ws.Range("M61").Value = -13164385
$ws.Range("N61").Value = -18354.111
$ws.Range("H74").Value = 46876424
$ws.Range("I74").Value = 53572770
$ws.Range("K74").Value = 53572770
$ws.Range("M74").Value = -53571896
$ws.Range("H77").Value = 46876424
$ws.Range("I77").Value = 53572770
$ws.Range("K77").Value = 267863850
$ws.Range("M77").Value = -267859482
$ws.Range("H97").Value = 609.5
$ws.Range("I97").Value = 609.5
$ws.Range("K97").Value = 609.5
$ws.Range("M97").Value = -113.5
$ws.Range("H116").Value = 791034.5
$ws.Range("I116").Value = 1086778.9
$ws.Range("J116").Value = 72798.14
$ws.Range("K116").Value = 1086778.9
$ws.Range("L116").Value = 72798.14
$ws.Range("M116").Value = -1084484.9
$ws.Range("N116").Value = -77386.14
$ws.Range("H132").Value = 23841.828
$ws.Range("I132").Value = 27978.762
$ws.Range("J132").Value = 12982.375
$ws.Range("K132").Value = 83936.28599999999
$ws.Range("L132").Value = 38947.125
$ws.Range("M132").Value = -81406.28599999999
$ws.Range("N132").Value = -44007.125
$ws.Range("H136").Value = 10647150
$ws.Range("I136").Value = 13164597
$ws.Range("J136").Value = 17930.111
$ws.Range("K136").Value = 39493791
$ws.Range("L136").Value = 53790.333
$ws.Range("M136").Value = -39491241
$ws.Range("N136").Value = -58890.333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 791034.5
$ws.Range("I3").Value = 1086778.9
$ws.Range("J3").Value = 72798.14
$ws.Range("K3").Value = 1086778.9
$ws.Range("L3").Value = 72798.14
$ws.Range("M3").Value = -1086664.9
$ws.Range("N3").Value = -73026.14
$ws.Range("H86").Value = 45456916
$ws.Range("I86").Value = 2169.6155
$ws.Range("K86").Value = 2169.6155
$ws.Range("M86").Value = -1046.6155
$ws.Range("H89").Value = 45456916
$ws.Range("I89").Value = 2169.6155
$ws.Range("K89").Value = 10848.0775
$ws.Range("M89").Value = -5232.077499999999
$ws.Range("H105").Value = 1616.7368
$ws.Range("I105").Value = 1110.6471
$ws.Range("K105").Value = 1110.6471
$ws.Range("M105").Value = 636.3529000000001
$ws.Range("H107").Value = 22728810
$ws.Range("I107").Value = 45456548
$ws.Range("J107").Value = 1073.091
$ws.Range("K107").Value = 45456548
$ws.Range("L107").Value = 1073.091
$ws.Range("M107").Value = -45454628
$ws.Range("N107").Value = -4913.091
$ws.Range("H140").Value = 181505.42
$ws.Range("J140").Value = 181505.42
$ws.Range("L140").Value = 181505.42
$ws.Range("N140").Value = -191865.42

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3209.75
$ws.Range("I31").Value = 2983.0833
$ws.Range("J31").Value = 5249.75
$ws.Range("K31").Value = 2983.0833
$ws.Range("L31").Value = 5249.75
$ws.Range("M31").Value = -2688.0833
$ws.Range("N31").Value = -5839.75
$ws.Range("H34").Value = 3209.75
$ws.Range("I34").Value = 2983.0833
$ws.Range("J34").Value = 5249.75
$ws.Range("K34").Value = 2983.0833
$ws.Range("L34").Value = 5249.75
$ws.Range("M34").Value = -2781.0833
$ws.Range("N34").Value = -5653.75
$ws.Range("H122").Value = 3666352.5
$ws.Range("I122").Value = 6414169
$ws.Range("K122").Value = 19242507
$ws.Range("M122").Value = -19240057
$ws.Range("H132").Value = 1710.4615
$ws.Range("I132").Value = 1137.3334
$ws.Range("K132").Value = 3412.0002
$ws.Range("M132").Value = -882.0001999999999
$ws.Range("H134").Value = 3396.8235
$ws.Range("I134").Value = 3096.077
$ws.Range("J134").Value = 4374.25
$ws.Range("K134").Value = 9288.231
$ws.Range("L134").Value = 13122.75
$ws.Range("M134").Value = -6753.231
$ws.Range("N134").Value = -18192.75
$ws.Range("H141").Value = 82492.17999999999
$ws.Range("J141").Value = 87891.2
$ws.Range("L141").Value = 87891.2
$ws.Range("N141").Value = -98251.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35326700
$ws.Range("I4").Value = 1040297.4
$ws.Range("K4").Value = 3120892.2
$ws.Range("M4").Value = -3120780.2
$ws.Range("H34").Value = 1720127.9
$ws.Range("J34").Value = 4900
$ws.Range("L34").Value = 14700
$ws.Range("N34").Value = -14868
$ws.Range("H39").Value = 4980
$ws.Range("J39").Value = 4991.6665
$ws.Range("L39").Value = 14974.9995
$ws.Range("N39").Value = -15562.9995
$ws.Range("H55").Value = 3606.3333
$ws.Range("J55").Value = 4947.125
$ws.Range("L55").Value = 14841.375
$ws.Range("N55").Value = -15195.375
$ws.Range("H122").Value = 476.7143
$ws.Range("J122").Value = 506
$ws.Range("L122").Value = 4554
$ws.Range("N122").Value = -9454

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 619723.3
$ws.Range("I80").Value = 1041070.7
$ws.Range("J80").Value = 6854.4546
$ws.Range("K80").Value = 1041070.7
$ws.Range("L80").Value = 6854.4546
$ws.Range("M80").Value = -1040072.7
$ws.Range("N80").Value = -8850.454600000001
$ws.Range("H83").Value = 619723.3
$ws.Range("I83").Value = 1041070.7
$ws.Range("J83").Value = 6854.4546
$ws.Range("K83").Value = 5205353.5
$ws.Range("L83").Value = 34272.273
$ws.Range("M83").Value = -5200361.5
$ws.Range("N83").Value = -44256.273
$ws.Range("H97").Value = 370.85715
$ws.Range("I97").Value = 411.6154
$ws.Range("K97").Value = 411.6154
$ws.Range("M97").Value = 84.38459999999998
$ws.Range("H132").Value = 6235.433
$ws.Range("I132").Value = 5673.2856
$ws.Range("K132").Value = 17019.8568
$ws.Range("M132").Value = -14489.8568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4979.85
$ws.Range("I46").Value = 1687.5
$ws.Range("J46").Value = 5345.6665
$ws.Range("K46").Value = 1687.5
$ws.Range("L46").Value = 5345.6665
$ws.Range("M46").Value = -1499.5
$ws.Range("N46").Value = -5721.6665
$ws.Range("H132").Value = 4941.7856
$ws.Range("I132").Value = 4325.85
$ws.Range("K132").Value = 12977.55
$ws.Range("M132").Value = -10447.55

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1607550.9
$ws.Range("I81").Value = 2088382.8
$ws.Range("K81").Value = 4176765.6
$ws.Range("M81").Value = -4175704.6
$ws.Range("H84").Value = 1607550.9
$ws.Range("I84").Value = 2088382.8
$ws.Range("K84").Value = 20883828
$ws.Range("M84").Value = -20878524
$ws.Range("H123").Value = 74990
$ws.Range("J123").Value = 74990
$ws.Range("L123").Value = 74990
$ws.Range("N123").Value = -84790
$ws.Range("H132").Value = 32685828
$ws.Range("I132").Value = 4275315
$ws.Range("K132").Value = 12825945
$ws.Range("M132").Value = -12823415
